$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.27464766666667
$ws.Range("H2").Value = 30.823943
$ws.Range("I2").Value = 0.9550701914680073
$ws.Range("J2").Value = 0.9550701914680074
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1795736666666667
$ws.Range("N2").Value = 0.538721
$ws.Range("O2").Value = 0.06721938218475064
$ws.Range("P2").Value = 0.06721938218475064
$ws.Range("Q2").Value = 1.845056155211444
$ws.Range("R2").Value = 16.605505396903
$ws.Range("S2").Value = 0.06419922821355095
$ws.Range("T2").Value = 0.06419922821355097
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.27464766666667
$ws.Range("H3").Value = 30.823943
$ws.Range("I3").Value = 0.9550701914680073
$ws.Range("J3").Value = 0.9550701914680074
$ws.Range("O3").Value = 0.3682385515018647
$ws.Range("P3").Value = 0.3682385515018648
$ws.Range("Q3").Value = 10.10751339795556
$ws.Range("R3").Value = 90.9676205816
$ws.Range("S3").Value = 0.3516936638887876
$ws.Range("T3").Value = 0.3516936638887877
$ws.Range("A4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.27464766666667
$ws.Range("H4").Value = 30.823943
$ws.Range("I4").Value = 0.9550701914680073
$ws.Range("J4").Value = 0.9550701914680074
$ws.Range("M4").Value = 1.481553
$ws.Range("N4").Value = 4.444659
$ws.Range("O4").Value = 0.5545861995390778
$ws.Range("P4").Value = 0.554586199539078
$ws.Range("Q4").Value = 15.222435074493
$ws.Range("R4").Value = 137.001915670437
$ws.Range("S4").Value = 0.5296687477793016
$ws.Range("T4").Value = 0.5296687477793017
$ws.Range("A5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.27464766666667
$ws.Range("H5").Value = 30.823943
$ws.Range("I5").Value = 0.9550701914680073
$ws.Range("J5").Value = 0.9550701914680074
$ws.Range("M5").Value = 0.02659666666666667
$ws.Range("N5").Value = 0.07979
$ws.Range("O5").Value = 0.00995586677430665
$ws.Range("P5").Value = 0.009955866774306652
$ws.Range("Q5").Value = 0.2732713791077778
$ws.Range("R5").Value = 2.45944241197
$ws.Range("S5").Value = 0.009508551586367025
$ws.Range("T5").Value = 0.009508551586367028
$ws.Range("A6").Value = "MuSCs"
$ws.Range("G6").Value = 0.475652
$ws.Range("H6").Value = 1.426956
$ws.Range("I6").Value = 0.04421378342596928
$ws.Range("J6").Value = 0.04421378342596929
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1795736666666667
$ws.Range("N6").Value = 0.538721
$ws.Range("O6").Value = 0.06721938218475064
$ws.Range("P6").Value = 0.06721938218475064
$ws.Range("Q6").Value = 0.08541457369733332
$ws.Range("R6").Value = 0.768731163276
$ws.Range("S6").Value = 0.002972023205944023
$ws.Range("T6").Value = 0.002972023205944023
$ws.Range("A7").Value = "MuSCs"
$ws.Range("G7").Value = 0.475652
$ws.Range("H7").Value = 1.426956
$ws.Range("I7").Value = 0.04421378342596928
$ws.Range("J7").Value = 0.04421378342596929
$ws.Range("O7").Value = 0.3682385515018647
$ws.Range("P7").Value = 0.3682385515018648
$ws.Range("Q7").Value = 0.4679147274666666
$ws.Range("R7").Value = 4.2112325472
$ws.Range("S7").Value = 0.01628121956519608
$ws.Range("T7").Value = 0.01628121956519609
$ws.Range("A8").Value = "MuSCs"
$ws.Range("G8").Value = 0.475652
$ws.Range("H8").Value = 1.426956
$ws.Range("I8").Value = 0.04421378342596928
$ws.Range("J8").Value = 0.04421378342596929
$ws.Range("M8").Value = 1.481553
$ws.Range("N8").Value = 4.444659
$ws.Range("O8").Value = 0.5545861995390778
$ws.Range("P8").Value = 0.554586199539078
$ws.Range("Q8").Value = 0.7047036475559999
$ws.Range("R8").Value = 6.342332828003999
$ws.Range("S8").Value = 0.02452035411745217
$ws.Range("T8").Value = 0.02452035411745218
$ws.Range("A9").Value = "MuSCs"
$ws.Range("G9").Value = 0.475652
$ws.Range("H9").Value = 1.426956
$ws.Range("I9").Value = 0.04421378342596928
$ws.Range("J9").Value = 0.04421378342596929
$ws.Range("M9").Value = 0.02659666666666667
$ws.Range("N9").Value = 0.07979
$ws.Range("O9").Value = 0.00995586677430665
$ws.Range("P9").Value = 0.009955866774306652
$ws.Range("Q9").Value = 0.01265075769333333
$ws.Range("R9").Value = 0.11385681924
$ws.Range("S9").Value = 0.0004401865373769976
$ws.Range("T9").Value = 0.0004401865373769978
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.007703000000000001
$ws.Range("H10").Value = 0.023109
$ws.Range("I10").Value = 0.0007160251060233982
$ws.Range("J10").Value = 0.0007160251060233983
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1795736666666667
$ws.Range("N10").Value = 0.538721
$ws.Range("O10").Value = 0.06721938218475064
$ws.Range("P10").Value = 0.06721938218475064
$ws.Range("Q10").Value = 0.001383255954333333
$ws.Range("R10").Value = 0.012449303589
$ws.Range("S10").Value = 0.0000481307652556634
$ws.Range("T10").Value = 0.0000481307652556634
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.007703000000000001
$ws.Range("H11").Value = 0.023109
$ws.Range("I11").Value = 0.0007160251060233982
$ws.Range("J11").Value = 0.0007160251060233983
$ws.Range("O11").Value = 0.3682385515018647
$ws.Range("P11").Value = 0.3682385515018648
$ws.Range("Q11").Value = 0.007577697866666668
$ws.Range("R11").Value = 0.0681992808
$ws.Range("S11").Value = 0.0002636680478810253
$ws.Range("T11").Value = 0.0002636680478810254
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.007703000000000001
$ws.Range("H12").Value = 0.023109
$ws.Range("I12").Value = 0.0007160251060233982
$ws.Range("J12").Value = 0.0007160251060233983
$ws.Range("M12").Value = 1.481553
$ws.Range("N12").Value = 4.444659
$ws.Range("O12").Value = 0.5545861995390778
$ws.Range("P12").Value = 0.554586199539078
$ws.Range("Q12").Value = 0.011412402759
$ws.Range("R12").Value = 0.102711624831
$ws.Range("S12").Value = 0.0003970976423240817
$ws.Range("T12").Value = 0.0003970976423240818
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.007703000000000001
$ws.Range("H13").Value = 0.023109
$ws.Range("I13").Value = 0.0007160251060233982
$ws.Range("J13").Value = 0.0007160251060233983
$ws.Range("M13").Value = 0.02659666666666667
$ws.Range("N13").Value = 0.07979
$ws.Range("O13").Value = 0.00995586677430665
$ws.Range("P13").Value = 0.009955866774306652
$ws.Range("Q13").Value = 0.0002048741233333334
$ws.Range("R13").Value = 0.00184386711
$ws.Range("S13").Value = 0.000007128650562627747
$ws.Range("T13").Value = 0.00000712865056262775
